# Sync attendance_reports: swap order of "System" and the email address
# in the "Recorded By" column (G) wherever the cell text is exactly
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
